$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated cryptos list values. Price cells (column D) are
# numeric-looking text (e.g. "0.580", "64.946.48" using "." as a thousands
# separator) so a leading apostrophe is used to force Excel to keep them
# as text, matching the original inlineStr string cells, instead of letting
# Excel auto-convert them to numbers (which would drop formatting such as
# trailing zeros).

$ws.Range("D2").Value = "'64.946.48"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "'3.464.87"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'577.23"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").Value = "'162.14"
$ws.Range("E6").Value = "  +3.94%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'3.466.12"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "  +8.81%  "
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").Value = "'0.440"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'4.056.84"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("E14").Value = "  -2.89%  "
$ws.Range("D15").Value = "'0.0000195"
$ws.Range("E15").Value = "  +5.41%  "
$ws.Range("D16").Value = "'28.74"
$ws.Range("E16").Value = "  +5.66%  "
$ws.Range("D17").Value = "'64.899.42"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "'3.462.59"
$ws.Range("E18").Value = "  +1.42%  "
$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'14.31"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'390.16"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  -3.38%  "
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").Value = "'73.09"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("E26").Value = "  +19.68%  "
$ws.Range("D27").Value = "'9.53"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").Value = "'0.182"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  +8.84%  "
$ws.Range("E31").Value = "  +6.30%  "
$ws.Range("D32").Value = "'2.05"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'23.68"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").Value = "'6.54"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D36").Value = "'7.09"
$ws.Range("E36").Value = "  +5.63%  "
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("E38").Value = "  +2.27%  "
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").Value = "'0.0770"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'2.954.23"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'27.42"
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "'4.56"
$ws.Range("E43").Value = "  +5.51%  "
$ws.Range("D44").Value = "'42.86"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "'0.777"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'24.12"
$ws.Range("E47").Value = "  +7.16%  "
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "'2.20"
$ws.Range("E49").Value = "  +11.26%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'0.873"
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("E51").Value = "  +3.94%  "
